# commit before refactor excel writeout
# Adds a "notes" column (D) to the model/yaml mapping sheet, fills in
# missing "include" (C) flags, and fills in missing "yaml" (B) values
# with "native_iso_EU27" for several rows; also adds a couple of notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D, then the two distinct note values, in the
# order they were first introduced.
$ws.Range("D1").Value  = "notes"
$ws.Range("D14").Value = "R10_Europe"
$ws.Range("D25").Value = "R10_Europe"
$ws.Range("D30").Value = "possibly only EU28?"

# Fill in missing yaml (column B) values
$ws.Range("B5").Value  = "native_iso_EU27"
$ws.Range("B7").Value  = "native_iso_EU27"
$ws.Range("B15").Value = "native_iso_EU27"
$ws.Range("B26").Value = "native_iso_EU27"
$ws.Range("B30").Value = "native_iso_EU27"
$ws.Range("B32").Value = "native_iso_EU27"
$ws.Range("B33").Value = "native_iso_EU27"
$ws.Range("B34").Value = "native_iso_EU27"

# Fill in missing include (column C) values
$ws.Range("C28").Value = "yes"
$ws.Range("C29").Value = "yes"
$ws.Range("C35").Value = "yes"
$ws.Range("C36").Value = "yes"
$ws.Range("C37").Value = "yes"
$ws.Range("C38").Value = "yes"
$ws.Range("C39").Value = "yes"
$ws.Range("C40").Value = "yes"
$ws.Range("C41").Value = "yes"
$ws.Range("C42").Value = "yes"
$ws.Range("C43").Value = "yes"

# Move the active selection to match the saved view state
$null = $ws.Range("B30").Select()
